$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking Price values (avoid float coercion)
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

# Apply the updated Price / Volume(1h) values
$ws.Range('D2').Value = '26.510.31'
$ws.Range('E2').Value = '  +0.24%  '
$ws.Range('D3').Value = '1.839.68'
$ws.Range('E3').Value = '  -0.06%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '260.44'
$ws.Range('E5').Value = '  -0.10%  '
$ws.Range('E6').Value = '  -0.05%  '
$ws.Range('D7').Value = '0.5259'
$ws.Range('E7').Value = '  +0.94%  '
$ws.Range('D8').Value = '0.3184'
$ws.Range('E8').Value = '  -1.65%  '
$ws.Range('E9').Value = '  -0.16%  '
$ws.Range('D10').Value = '18.75'
$ws.Range('E10').Value = '  +0.73%  '
$ws.Range('E11').Value = '  +2.29%  '
$ws.Range('D12').Value = '0.07754'
$ws.Range('E12').Value = '  +0.86%  '
$ws.Range('D13').Value = '1.842.71'
$ws.Range('E13').Value = '  -0.09%  '
$ws.Range('D14').Value = '87.87'
$ws.Range('E14').Value = '  -0.65%  '
$ws.Range('D15').Value = '5.011'
$ws.Range('E15').Value = '  -0.25%  '
$ws.Range('E16').Value = '  -0.14%  '
$ws.Range('D17').Value = '13.84'
$ws.Range('E17').Value = '  -0.57%  '
$ws.Range('E18').Value = '  +0.00%  '
$ws.Range('D19').Value = '0.000007951'
$ws.Range('E19').Value = '  -0.03%  '
$ws.Range('D20').Value = '26.534.85'
$ws.Range('E20').Value = '  +0.26%  '
$ws.Range('D21').Value = '2.072.88'
$ws.Range('E21').Value = '  -0.05%  '
$ws.Range('D22').Value = '4.617'
$ws.Range('E22').Value = '  +1.16%  '
$ws.Range('D23').Value = '5.975'
$ws.Range('E23').Value = '  +0.57%  '
$ws.Range('D24').Value = '9.319'
$ws.Range('D25').Value = '142.00'
$ws.Range('E25').Value = '  -1.76%  '
$ws.Range('D26').Value = '2.200'
$ws.Range('E26').Value = '  -0.47%  '
$ws.Range('E27').Value = '  +1.47%  '
$ws.Range('D28').Value = '16.89'
$ws.Range('E28').Value = '  -0.51%  '
$ws.Range('D29').Value = '111.45'
$ws.Range('E29').Value = '  +0.15%  '
$ws.Range('E30').Value = '  +0.07%  '
$ws.Range('D31').Value = '0.08705'
$ws.Range('E31').Value = '  -0.31%  '
$ws.Range('D32').Value = '4.074'
$ws.Range('E32').Value = '  -1.14%  '
$ws.Range('D33').Value = '0.04881'
$ws.Range('E33').Value = '  +1.76%  '
$ws.Range('D34').Value = '0.7268'
$ws.Range('E34').Value = '  +4.16%  '
$ws.Range('D35').Value = '1.133'
$ws.Range('E35').Value = '  +1.05%  '
$ws.Range('E36').Value = '  +0.24%  '
$ws.Range('D37').Value = '3.093'
$ws.Range('E37').Value = '  +1.15%  '
$ws.Range('D38').Value = '2.235'
$ws.Range('E38').Value = '  +2.37%  '
$ws.Range('E39').Value = '  -0.07%  '
$ws.Range('D40').Value = '0.4809'
$ws.Range('E40').Value = '  -0.28%  '
$ws.Range('D41').Value = '0.8938'
$ws.Range('E41').Value = '  +0.53%  '
$ws.Range('D42').Value = '109.47'
$ws.Range('E42').Value = '  -1.36%  '
$ws.Range('E43').Value = '  -2.38%  '
$ws.Range('E44').Value = '  +0.00%  '
$ws.Range('D45').Value = '7.643'
$ws.Range('E45').Value = '  -0.12%  '
$ws.Range('D46').Value = '0.4178'
$ws.Range('E46').Value = '  +1.24%  '
$ws.Range('D47').Value = '8.995'
$ws.Range('E47').Value = '  -0.38%  '
$ws.Range('D48').Value = '0.05845'
$ws.Range('E48').Value = '  -0.18%  '
$ws.Range('D49').Value = '0.1233'
$ws.Range('E49').Value = '  +1.20%  '
$ws.Range('D50').Value = '34.88'
$ws.Range('E50').Value = '  -0.12%  '
$ws.Range('D51').Value = '0.8914'
$ws.Range('E51').Value = '  +0.98%  '

# Restore default styling (no visible/semantic format change intended)
$ws.Range('D5').Style = 'Normal'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
